# Test Cases.xlsx edit:
#   "Added remaining test cases & updated case IDs"
#
# The sheet had three blank rows (4:6) separating the first test case
# (CreateAccount_1, row 3) from the rest of the table (old rows 7:13).
# This edit:
#   1. Removes the blank rows 4:6, pulling the remaining test cases up.
#   2. Bumps the numeric suffix on each TestCaseId (e.g. CreateOrder_1 ->
#      CreateOrder_2, SystemTest_1 -> SystemTest_3, ... PasswordReset_6 ->
#      PasswordReset_8).
#   3. Appends two brand-new test cases (PromoCode_9, PromoCode_10) as
#      rows 11 and 12, copying the visual style used by the other
#      "full border" rows (same style family as row 3 / CreateAccount_1).
#   4. Restores the row heights Excel recalculated for the new layout.
#   5. Leaves the active selection on M5 (matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the 3 empty separator rows -----------------------------------
# Old layout: row3=CreateAccount_1, rows4-6=blank, row7=CreateOrder_1, ...
# After deleting 4:6, the old row7..row13 content becomes row4..row10.
$ws.Range("A4:A6").EntireRow.Delete()

# --- 2. Bump the TestCaseId suffix on each pre-existing case --------------
$ws.Range("A4").Value  = "CreateOrder_2"     # was CreateOrder_1
$ws.Range("A5").Value  = "SystemTest_3"      # was SystemTest_1
$ws.Range("A6").Value  = "Login Auth_4"      # was Login Auth_2
$ws.Range("A7").Value  = "SeatSelect_5"      # was SeatSelect_3
$ws.Range("A8").Value  = "Payment_6"         # was Payment_4
$ws.Range("A9").Value  = "Refund_7"          # was Refund_5
$ws.Range("A10").Value = "PasswordReset_8"   # was PasswordReset_6

# --- 3. Append the two new test cases (rows 11 & 12) -----------------------
# Row 3 uses the bordered/wrapped style family (s=2/3/3/4/3/3/5/2) that the
# two new rows also use, so clone its formatting onto the new rows first.
$ws.Range("A3:J3").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A3:J3").Copy()
$ws.Range("A12:J12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "PromoCode_9"
$ws.Range("B11").Value = "Payment Gateway"
$ws.Range("C11").Value = "P2"
$ws.Range("D11").Value = "Verify that valid promo codes correctly apply discounts correctly and invalid codes display an error."
$ws.Range("E11").Value = "User has selected tickets and reached the payment page"
$ws.Range("F11").Value = "1. On the payment page, enter a valid promo code.
2. Verify that the total cost updates to the correct amount after the discount is applied.
3. Enter an invalid promo code.
4. Try to complete the payment."
$ws.Range("G11").Value = "After entering the valid promo code, the total cust shows the correct amount after the discount is applied and the invalid promo code displays an error."
$ws.Range("H11").Value = "Not yet tested"
$ws.Range("I11").Value = "Not Executed"
$ws.Range("J11").Value = "James Shumate"

$ws.Range("A12").Value = "PromoCode_10"
$ws.Range("B12").Value = "User and Payment Gateway"
$ws.Range("C12").Value = "P1"
$ws.Range("D12").Value = "Verify integration between User and Payment modules, ensuring purchases are correctly displayed in the user's purchase history."
$ws.Range("E12").Value = "User is logged in to a user account."
$ws.Range("F12").Value = "1. Log in as a user.
2. Select a theater, then movie, then select seats.
3. Complete the purchase.
4. Navigate to the user's purchase history.
5. Confirm that the purchase appears with correct details."
$ws.Range("G12").Value = "The purchase correctly appears in the user's purchase history. The data is consistent between Payment and User modules."
$ws.Range("H12").Value = "Not yet tested"
$ws.Range("I12").Value = "Not Executed"
$ws.Range("J12").Value = "James Shumate"

# --- 4. Row heights (re-wrapped text needs new heights) --------------------
$ws.Rows.Item(1).RowHeight  = 21
$ws.Rows.Item(3).RowHeight  = 80.25
$ws.Rows.Item(4).RowHeight  = 114
$ws.Rows.Item(5).RowHeight  = 136.5
$ws.Rows.Item(6).RowHeight  = 91.5
$ws.Rows.Item(7).RowHeight  = 57.75
$ws.Rows.Item(8).RowHeight  = 91.5
$ws.Rows.Item(9).RowHeight  = 80.25
$ws.Rows.Item(10).RowHeight = 91.5
$ws.Rows.Item(11).RowHeight = 125.25
$ws.Rows.Item(12).RowHeight = 125.25

# --- 5. Restore the saved selection -----------------------------------------
$ws.Range("M5").Select()
